# Update the cryptos price/volume table with the latest scraped values.
# For the "Price" column (D) the values are stored as literal text (e.g. "36.662.50"
# is not a valid number). A leading apostrophe is prepended before assignment so
# Excel keeps them as exact text (preserving trailing zeros, multiple dots, etc.)
# instead of auto-coercing them into Double values via COM automation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '36.602.16'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = "'" + '1.962.70'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'" + '244.47'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('E6').Value = '  +1.51%  '
$ws.Range('D7').Value = "'" + '60.49'
$ws.Range('E7').Value = '  +7.33%  '
$ws.Range('E9').Value = '  +5.14%  '
$ws.Range('D10').Value = "'" + '0.0796'
$ws.Range('E10').Value = '  -4.40%  '
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').Value = "'" + '14.30'
$ws.Range('E12').Value = '  +7.59%  '
$ws.Range('D13').Value = "'" + '0.841'
$ws.Range('E13').Value = '  +5.35%  '
$ws.Range('D14').Value = "'" + '21.83'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = "'" + '2.251.58'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = "'" + '5.31'
$ws.Range('E16').Value = '  +3.94%  '
$ws.Range('D17').Value = "'" + '1.958.26'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = "'" + '36.585.61'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').Value = "'" + '70.03'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').Value = "'" + '0.0₃0855'
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = "'" + '230.63'
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('D22').Value = "'" + '5.10'
$ws.Range('E22').Value = '  +3.75%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +6.48%  '
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').Value = "'" + '0.145'
$ws.Range('E26').Value = '  +11.60%  '
$ws.Range('D27').Value = "'" + '9.23'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('D28').Value = "'" + '160.91'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').Value = "'" + '19.46'
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').Value = "'" + '1.23'
$ws.Range('E30').Value = '  +12.31%  '
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('D32').Value = "'" + '4.77'
$ws.Range('E32').Value = '  +6.00%  '
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = "'" + '4.47'
$ws.Range('E34').Value = '  +8.13%  '
$ws.Range('D35').Value = "'" + '3.62'
$ws.Range('E35').Value = '  +22.08%  '
$ws.Range('E36').Value = '  +6.48%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = "'" + '1.77'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').Value = "'" + '5.58'
$ws.Range('E39').Value = '  -5.67%  '
$ws.Range('D40').Value = "'" + '0.0985'
$ws.Range('E41').Value = '  +1.40%  '
$ws.Range('E42').Value = '  +3.77%  '
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').Value = "'" + '16.24'
$ws.Range('E44').Value = '  +4.35%  '
$ws.Range('D45').Value = "'" + '1.367.50'
$ws.Range('E45').Value = '  +3.46%  '
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').Value = "'" + '88.68'
$ws.Range('E47').Value = '  +4.86%  '
$ws.Range('D48').Value = "'" + '7.17'
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('D50').Value = "'" + '44.37'
$ws.Range('E50').Value = '  +3.10%  '
$ws.Range('D51').Value = "'" + '1.84'
$ws.Range('E51').Value = '  +6.02%  '
